$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.251.29'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.36%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.906.62'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.11%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '327.59'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4658'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.65%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3959'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.34%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.70'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.75%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07966'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.41%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.004'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.08%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.36'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.907.18'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.137'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.40%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.789'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.89%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06955'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.10%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '88.95'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.91%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.002'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001013'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.31%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.18'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.29%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.003'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '29.251.65'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.35%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.364'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.11'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.71%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.109.13'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.37%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.058'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.32%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.79'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.53'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.64%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.897'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.99%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.007'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '119.50'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.35%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09430'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.29%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9255'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.355'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.66%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.350'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.40%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.258'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.89%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05860'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.39%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.170'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.89%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02112'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.16%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.987'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.83%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5763'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.67%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1814'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.87%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '10.03'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.92%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.02'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.89%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5437'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.91%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.232'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.33%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.07098'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.03%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.888'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.16%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.589'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.35%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '112.06'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.072'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.79%  '
